# Meniere's Disease Map title textbox: shrink font 28pt -> 24pt and, since the
# textbox uses <a:spAutoFit/>, resize/reposition it (and the enclosing group's
# bounding box) to match how PowerPoint auto-fits the box to the smaller text.
#
# EMU -> point conversion uses 12700 EMU per point. A tiny epsilon is added
# before assignment to counter float32 round-trip truncation in this COM
# runtime so the saved EMU values land exactly on the intended integers.

$EMU_PER_PT = 12700
$eps = 0.000002

function ToPt([double]$emu) {
    return ($emu / $EMU_PER_PT) + $eps
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Group 4" is the only top-level shape on the slide; it contains the title
# rectangle ("Rectangulo 18") and the SmartArt graphic frame ("Diagrama 1").
$grp = $s.Shapes.Item(1)
$titleBox = $grp.GroupItems.Item(1)

# 1) Shrink the title font from 28pt to 24pt.
$titleBox.TextFrame.TextRange.Font.Size = 24

# 2) Resize/reposition the title textbox to its new auto-fit bounds.
$titleBox.Left = ToPt(2897445)
$titleBox.Top = ToPt(950181)
$titleBox.Width = ToPt(3156826)
$titleBox.Height = ToPt(461665)

# 3) The group's own bounding box (off/ext/chOff/chExt) must now equal the
# union of its children's bounds. Ungrouping and regrouping recomputes this
# bounding box from the (now-updated) children, then restore the original
# group name.
$ungrouped = $grp.Ungroup()
$regrouped = $ungrouped.Group()
$regrouped.Name = "Group 4"
